$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency Price (D) / Volume(1h) (E) figures.
# Some new Price strings are valid-looking numbers (e.g. "154.53");
# prefix those with a leading apostrophe so Excel stores them as
# text, matching the original inlineStr string cells, rather than
# silently converting them to the Number type.

$ws.Range("D2").Value = "65.924.82"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.669.78"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'599.01"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'158.00"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").Value = "'0.651"
$ws.Range("E7").Value = "  +4.48%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.127"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").Value = "'5.85"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "'29.16"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "'0.0000195"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "3.147.45"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "65.785.25"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "2.668.92"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "'12.65"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").Value = "'4.81"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "'350.71"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "'7.47"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'69.80"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("E24").Value = "  +10.73%  "
$ws.Range("D25").Value = "'0.0000113"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Value = "'9.64"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "'1.64"
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("D28").Value = "'568.38"
$ws.Range("E28").Value = "  +6.55%  "
$ws.Range("D29").Value = "'8.25"
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").Value = "'0.164"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "'1.83"
$ws.Range("E33").Value = "  +3.65%  "
$ws.Range("D34").Value = "'6.70"
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("D35").Value = "'5.56"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("D36").Value = "'0.424"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "'20.61"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'1.96"
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").Value = "'154.53"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("D41").Value = "'161.31"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("D42").Value = "'4.11"
$ws.Range("D43").Value = "'0.0621"
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("D44").Value = "'2.33"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'23.06"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("D46").Value = "'0.644"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "'0.0256"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "'19.89"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "0.0₆0244"
$ws.Range("E50").Value = "  -5.57%  "
$ws.Range("D51").Value = "'0.815"
$ws.Range("E51").Value = "  -0.74%  "
